$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 373.1111
$ws.Cells.Item(2, 9).Value = 373.1111
$ws.Cells.Item(2, 11).Value = 373.1111
$ws.Cells.Item(2, 13).Value = -260.1111
$ws.Cells.Item(40, 8).Value = 2893.2942
$ws.Cells.Item(40, 9).Value = 2497.8
$ws.Cells.Item(40, 11).Value = 2497.8
$ws.Cells.Item(40, 13).Value = -2322.8
$ws.Cells.Item(58, 8).Value = 221.09091
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 12).Value = 3000
$ws.Cells.Item(58, 14).Value = -3300
$ws.Cells.Item(64, 8).Value = 6495.8335
$ws.Cells.Item(64, 9).Value = 4132
$ws.Cells.Item(64, 11).Value = 4132
$ws.Cells.Item(64, 13).Value = -3884
$ws.Cells.Item(67, 8).Value = 6495.8335
$ws.Cells.Item(67, 9).Value = 4132
$ws.Cells.Item(67, 11).Value = 4132
$ws.Cells.Item(67, 13).Value = -3274
$ws.Cells.Item(70, 8).Value = 1461366.9
$ws.Cells.Item(70, 9).Value = 3402189.8
$ws.Cells.Item(70, 10).Value = 5749.75
$ws.Cells.Item(70, 11).Value = 10206569.4
$ws.Cells.Item(70, 12).Value = 17249.25
$ws.Cells.Item(70, 13).Value = -10206299.4
$ws.Cells.Item(70, 14).Value = -17789.25
$ws.Cells.Item(73, 8).Value = 1461366.9
$ws.Cells.Item(73, 9).Value = 3402189.8
$ws.Cells.Item(73, 10).Value = 5749.75
$ws.Cells.Item(73, 11).Value = 10206569.4
$ws.Cells.Item(73, 12).Value = 17249.25
$ws.Cells.Item(73, 13).Value = -10205633.4
$ws.Cells.Item(73, 14).Value = -19121.25
$ws.Cells.Item(107, 8).Value = 384.25
$ws.Cells.Item(107, 9).Value = 384.25
$ws.Cells.Item(107, 11).Value = 384.25
$ws.Cells.Item(107, 13).Value = 1535.75
$ws.Cells.Item(111, 8).Value = 30750.637
$ws.Cells.Item(111, 9).Value = 1020.375
$ws.Cells.Item(111, 10).Value = 110031.336
$ws.Cells.Item(111, 11).Value = 3061.125
$ws.Cells.Item(111, 12).Value = 330094.008
$ws.Cells.Item(111, 13).Value = 5.875
$ws.Cells.Item(111, 14).Value = -336228.008
$ws.Cells.Item(133, 8).Value = 96649
$ws.Cells.Item(133, 10).Value = 99999
$ws.Cells.Item(133, 12).Value = 99999
$ws.Cells.Item(133, 14).Value = -110119
$ws.Cells.Item(138, 8).Value = 2126.85
$ws.Cells.Item(138, 10).Value = 2556.3572
$ws.Cells.Item(138, 12).Value = 7669.071599999999
$ws.Cells.Item(138, 14).Value = -17949.0716

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9728.831
$ws.Cells.Item(32, 9).Value = 6345.1045
$ws.Cells.Item(32, 11).Value = 6345.1045
$ws.Cells.Item(32, 13).Value = -6058.1045
$ws.Cells.Item(45, 8).Value = 328428.6
$ws.Cells.Item(45, 10).Value = 1774.4166
$ws.Cells.Item(45, 12).Value = 1774.4166
$ws.Cells.Item(45, 14).Value = -2528.4166
$ws.Cells.Item(97, 8).Value = 1894.5555
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 4834022
$ws.Cells.Item(102, 9).Value = 5708807.5
$ws.Cells.Item(102, 10).Value = 168499.67
$ws.Cells.Item(102, 11).Value = 5708807.5
$ws.Cells.Item(102, 12).Value = 168499.67
$ws.Cells.Item(102, 13).Value = -5707185.5
$ws.Cells.Item(102, 14).Value = -171743.67
$ws.Cells.Item(110, 8).Value = 45456144
$ws.Cells.Item(110, 9).Value = 52632852
$ws.Cells.Item(110, 11).Value = 52632852
$ws.Cells.Item(110, 13).Value = -52630807
$ws.Cells.Item(122, 8).Value = 4292.5293
$ws.Cells.Item(122, 9).Value = 3839.8696
$ws.Cells.Item(122, 11).Value = 11519.6088
$ws.Cells.Item(122, 13).Value = -9069.6088

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 10151.846
$ws.Cells.Item(20, 9).Value = 3203.1428
$ws.Cells.Item(20, 11).Value = 3203.1428
$ws.Cells.Item(20, 13).Value = -2956.1428
$ws.Cells.Item(22, 8).Value = 171
$ws.Cells.Item(22, 9).Value = 88.75
$ws.Cells.Item(22, 11).Value = 88.75
$ws.Cells.Item(22, 13).Value = 84.25
$ws.Cells.Item(105, 8).Value = 2169.4583
$ws.Cells.Item(105, 9).Value = 1964.3529
$ws.Cells.Item(105, 10).Value = 2667.5715
$ws.Cells.Item(105, 11).Value = 1964.3529
$ws.Cells.Item(105, 12).Value = 2667.5715
$ws.Cells.Item(105, 13).Value = -217.3529000000001
$ws.Cells.Item(105, 14).Value = -6161.5715
$ws.Cells.Item(107, 8).Value = 31252272
$ws.Cells.Item(107, 9).Value = 1966.7273
$ws.Cells.Item(107, 11).Value = 1966.7273
$ws.Cells.Item(107, 13).Value = -46.72730000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(21, 8).Value = 9999
$ws.Cells.Item(21, 9).Value = 9999
$ws.Cells.Item(21, 11).Value = 9999
$ws.Cells.Item(21, 13).Value = -9764
$ws.Cells.Item(31, 8).Value = 3183.7097
$ws.Cells.Item(31, 9).Value = 2410.524
$ws.Cells.Item(31, 10).Value = 4807.4
$ws.Cells.Item(31, 11).Value = 2410.524
$ws.Cells.Item(31, 12).Value = 4807.4
$ws.Cells.Item(31, 13).Value = -2115.524
$ws.Cells.Item(31, 14).Value = -5397.4
$ws.Cells.Item(34, 8).Value = 3183.7097
$ws.Cells.Item(34, 9).Value = 2410.524
$ws.Cells.Item(34, 10).Value = 4807.4
$ws.Cells.Item(34, 11).Value = 2410.524
$ws.Cells.Item(34, 12).Value = 4807.4
$ws.Cells.Item(34, 13).Value = -2208.524
$ws.Cells.Item(34, 14).Value = -5211.4
$ws.Cells.Item(62, 8).Value = 118842.38
$ws.Cells.Item(62, 9).Value = 4626.25
$ws.Cells.Item(62, 10).Value = 169605.11
$ws.Cells.Item(62, 11).Value = 4626.25
$ws.Cells.Item(62, 12).Value = 169605.11
$ws.Cells.Item(62, 13).Value = -4002.25
$ws.Cells.Item(62, 14).Value = -170853.11
$ws.Cells.Item(65, 8).Value = 118842.38
$ws.Cells.Item(65, 9).Value = 4626.25
$ws.Cells.Item(65, 10).Value = 169605.11
$ws.Cells.Item(65, 11).Value = 23131.25
$ws.Cells.Item(65, 12).Value = 848025.5499999999
$ws.Cells.Item(65, 13).Value = -20011.25
$ws.Cells.Item(65, 14).Value = -854265.5499999999
$ws.Cells.Item(141, 8).Value = 215976.14
$ws.Cells.Item(141, 10).Value = 215976.14
$ws.Cells.Item(141, 12).Value = 215976.14
$ws.Cells.Item(141, 14).Value = -226336.14

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 274.68
$ws.Cells.Item(4, 9).Value = 182.93828
$ws.Cells.Item(4, 10).Value = 665.7895
$ws.Cells.Item(4, 11).Value = 548.81484
$ws.Cells.Item(4, 12).Value = 1997.3685
$ws.Cells.Item(4, 13).Value = -436.81484
$ws.Cells.Item(4, 14).Value = -2221.3685
$ws.Cells.Item(5, 8).Value = 900
$ws.Cells.Item(5, 9).Value = 800
$ws.Cells.Item(5, 11).Value = 2400
$ws.Cells.Item(5, 13).Value = -2288
$ws.Cells.Item(34, 8).Value = 2584.4614
$ws.Cells.Item(34, 10).Value = 2936.182
$ws.Cells.Item(34, 12).Value = 8808.545999999998
$ws.Cells.Item(34, 14).Value = -8976.545999999998
$ws.Cells.Item(131, 8).Value = 10163.5
$ws.Cells.Item(131, 9).Value = 1694.5
$ws.Cells.Item(131, 10).Value = 15456.625
$ws.Cells.Item(131, 11).Value = 5083.5
$ws.Cells.Item(131, 12).Value = 46369.875
$ws.Cells.Item(131, 13).Value = -43.5
$ws.Cells.Item(131, 14).Value = -56449.875
$ws.Cells.Item(135, 8).Value = 900
$ws.Cells.Item(135, 9).Value = 800
$ws.Cells.Item(135, 11).Value = 7200
$ws.Cells.Item(135, 13).Value = -4665

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 20229.334
$ws.Cells.Item(95, 10).Value = 20229.334
$ws.Cells.Item(95, 12).Value = 20229.334
$ws.Cells.Item(95, 14).Value = -25721.334
$ws.Cells.Item(113, 8).Value = 3780.8235
$ws.Cells.Item(113, 10).Value = 4150.4
$ws.Cells.Item(113, 12).Value = 4150.4
$ws.Cells.Item(113, 14).Value = -8490.4
$ws.Cells.Item(122, 8).Value = 208184.27
$ws.Cells.Item(122, 9).Value = 280038.16
$ws.Cells.Item(122, 10).Value = 2887.4285
$ws.Cells.Item(122, 11).Value = 840114.48
$ws.Cells.Item(122, 12).Value = 8662.2855
$ws.Cells.Item(122, 13).Value = -837664.48
$ws.Cells.Item(122, 14).Value = -13562.2855

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 2500
$ws.Cells.Item(25, 9).Value = 2500
$ws.Cells.Item(25, 11).Value = 2500
$ws.Cells.Item(25, 13).Value = -2270
$ws.Cells.Item(55, 8).Value = 900.5
$ws.Cells.Item(55, 9).Value = 545
$ws.Cells.Item(55, 11).Value = 545
$ws.Cells.Item(55, 13).Value = -372
$ws.Cells.Item(61, 8).Value = 862.2
$ws.Cells.Item(61, 9).Value = 832.75
$ws.Cells.Item(61, 11).Value = 832.75
$ws.Cells.Item(61, 13).Value = -630.75
$ws.Cells.Item(113, 8).Value = 862.2
$ws.Cells.Item(113, 9).Value = 832.75
$ws.Cells.Item(113, 11).Value = 832.75
$ws.Cells.Item(113, 13).Value = 1337.25
$ws.Cells.Item(133, 8).Value = 79762.60000000001
$ws.Cells.Item(133, 10).Value = 79762.60000000001
$ws.Cells.Item(133, 12).Value = 79762.60000000001
$ws.Cells.Item(133, 14).Value = -84822.60000000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 947.96
$ws.Cells.Item(107, 9).Value = 813
$ws.Cells.Item(107, 11).Value = 2439
$ws.Cells.Item(107, 13).Value = -519
$ws.Cells.Item(113, 8).Value = 14133.6
$ws.Cells.Item(113, 9).Value = 26209.5
$ws.Cells.Item(113, 10).Value = 6083
$ws.Cells.Item(113, 11).Value = 78628.5
$ws.Cells.Item(113, 12).Value = 18249
$ws.Cells.Item(113, 13).Value = -76458.5
$ws.Cells.Item(113, 14).Value = -22589

Write-Host "All edits applied"